$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen columns J:AA (content was added requiring more room) ---
$ws.Range("J1").EntireColumn.ColumnWidth = 21.3307
$ws.Range("K1").EntireColumn.ColumnWidth = 31.3307
$ws.Range("L1").EntireColumn.ColumnWidth = 26.8307
$ws.Range("M1").EntireColumn.ColumnWidth = 32.4987
$ws.Range("N1").EntireColumn.ColumnWidth = 28.8307
$ws.Range("O1").EntireColumn.ColumnWidth = 21.8307
$ws.Range("P1").EntireColumn.ColumnWidth = 35.8307
$ws.Range("Q1").EntireColumn.ColumnWidth = 28.8307
$ws.Range("R1").EntireColumn.ColumnWidth = 24.6667
$ws.Range("S1").EntireColumn.ColumnWidth = 25.3307
$ws.Range("T1").EntireColumn.ColumnWidth = 23.1667
$ws.Range("U1").EntireColumn.ColumnWidth = 21.9987
$ws.Range("V1").EntireColumn.ColumnWidth = 19.3307
$ws.Range("W1").EntireColumn.ColumnWidth = 17.6667
$ws.Range("X1").EntireColumn.ColumnWidth = 19.3307
$ws.Range("Y1").EntireColumn.ColumnWidth = 18.4987
$ws.Range("Z1").EntireColumn.ColumnWidth = 16.8307
$ws.Range("AA1").EntireColumn.ColumnWidth = 16.6667

# --- Replace the static AD (total score) values with a live SUM formula ---
# Mirrors how it was actually entered: AD2 on its own, then AD3:AD66 filled
# down together, then AD67:AD130 filled down together, then AD131 on its own.
$ws.Range("AD2").Formula = "=SUM(I2:AC2)"
$ws.Range("AD3:AD66").Formula = "=SUM(I3:AC3)"
$ws.Range("AD67:AD130").Formula = "=SUM(I67:AC67)"
$ws.Range("AD131").Formula = "=SUM(I131:AC131)"

# --- Update the active selection ---
[void]$ws.Range("AD15").Select()
